$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal": update the figures for the current (last) monthly row ---
$wsMensal = $wb.Worksheets.Item("Mensal")

$wsMensal.Cells.Item(14, 1).Value = 44066
$wsMensal.Cells.Item(14, 2).Value = 157.98
$wsMensal.Cells.Item(14, 4).Value = -9.300000000000001

# --- Sheet "Diario": append the newly reported daily rows (384-390) ---
$wsDiario = $wb.Worksheets.Item("Diario")

$dados = @(
    @(44060, 161.74, 174.17, -7.14),
    @(44061, 164.4,  174.17, -5.61),
    @(44062, 198.93, 174.17, 14.21),
    @(44063, 187.89, 174.17, 7.87),
    @(44064, 169,    174.17, -2.97),
    @(44065, 163.08, 174.17, -6.37),
    @(44066, 160.41, 174.17, -7.9)
)

$lastRow = 383
for ($i = 0; $i -lt $dados.Length; $i++) {
    $r = $lastRow + $i + 1
    $linha = $dados[$i]

    # Copy formatting (incl. the date number format/border/font on column A)
    # from the previous row down into the new row before writing values.
    $wsDiario.Range("A$($r - 1):D$($r - 1)").Copy($wsDiario.Range("A$($r):D$($r)"))

    $wsDiario.Cells.Item($r, 1).Value = $linha[0]
    $wsDiario.Cells.Item($r, 2).Value = $linha[1]
    $wsDiario.Cells.Item($r, 3).Value = $linha[2]
    $wsDiario.Cells.Item($r, 4).Value = $linha[3]
}
